$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column "Q" = "Дата вывода из эксплуатации" (archiveDate) after the
# existing "Автоматический" column (P), mirroring its formatting (style 7),
# and fill in the sample archive date only for the first data row, same as
# the other sample date columns.

# 1) Copy the formatting (fill/border/number format) from column P onto the
#    new column Q for header row + all 10 data rows so the new cells get the
#    same visual style (s="7") as their P-column neighbours.
$ws.Range("P1:P11").Copy() | Out-Null
$ws.Range("Q1:Q11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# 2) Header cell
$ws.Range("Q1").NumberFormat = "@"
$ws.Range("Q1").Value = "Дата вывода из эксплуатации"

# 3) Sample value - only row 2 gets a value in the example file
$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "2022-01-25"

# 4) Column width for the new column (closest achievable width to the
#    original authored 24.1719 given this engine's column-width rounding)
$ws.Columns("Q").ColumnWidth = 24.0833
